# Update test data for number input of LOINC and SNOMED code.
$wb = $excel.ActiveWorkbook

$wsDescriptor = $wb.Worksheets.Item("descriptor")
$wsProperty   = $wb.Worksheets.Item("property")

# "WALK" row in the descriptor sheet: SNOMED_code is now entered as a
# plain number (100) instead of the placeholder text "S001".
$wsDescriptor.Range("I2").Value = 100

# "STEPS" row in the property sheet: both LOINC_code and SNOMED_code are
# now entered as plain numbers (100) instead of the placeholder text
# "L001" / "S001".
$wsProperty.Range("B2").Value = 100
$wsProperty.Range("E2").Value = 100

# Reflect the author's final selection/view state: the property sheet was
# zoomed out and left with E3 selected (not the active tab)...
$wsProperty.Range("E3").Select()
$excel.ActiveWindow.Zoom = 180

# ...while the descriptor sheet stayed the active tab, selection moved to H9.
$wsDescriptor.Activate()
$wsDescriptor.Range("H9").Select()
